# Update the worksheet date and regenerate the per-problem operands while
# leaving every other cell (e.g. "553×4=") untouched. Each "old" string is
# unique across the document, so a simple Find/Replace-All per pair is safe.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-24 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-25 Monday", 2) | Out-Null
$d.Content.Find.Execute("544×4=", $true, $false, $false, $false, $false, $true, 1, $false, "271×5=", 2) | Out-Null
$d.Content.Find.Execute("984×6=", $true, $false, $false, $false, $false, $true, 1, $false, "924×9=", 2) | Out-Null
$d.Content.Find.Execute("293×9=", $true, $false, $false, $false, $false, $true, 1, $false, "418×3=", 2) | Out-Null
$d.Content.Find.Execute("861×3=", $true, $false, $false, $false, $false, $true, 1, $false, "571×6=", 2) | Out-Null
$d.Content.Find.Execute("139×4=", $true, $false, $false, $false, $false, $true, 1, $false, "155×9=", 2) | Out-Null
$d.Content.Find.Execute("657×8=", $true, $false, $false, $false, $false, $true, 1, $false, "969×3=", 2) | Out-Null
$d.Content.Find.Execute("307×9=", $true, $false, $false, $false, $false, $true, 1, $false, "637×2=", 2) | Out-Null
$d.Content.Find.Execute("281×5=", $true, $false, $false, $false, $false, $true, 1, $false, "378×7=", 2) | Out-Null
$d.Content.Find.Execute("845×4=", $true, $false, $false, $false, $false, $true, 1, $false, "839×9=", 2) | Out-Null
$d.Content.Find.Execute("160×4=", $true, $false, $false, $false, $false, $true, 1, $false, "110×7=", 2) | Out-Null
$d.Content.Find.Execute("500×8=", $true, $false, $false, $false, $false, $true, 1, $false, "586×9=", 2) | Out-Null
$d.Content.Find.Execute("446×9=", $true, $false, $false, $false, $false, $true, 1, $false, "513×3=", 2) | Out-Null
$d.Content.Find.Execute("420×4=", $true, $false, $false, $false, $false, $true, 1, $false, "693×6=", 2) | Out-Null
$d.Content.Find.Execute("900×4=", $true, $false, $false, $false, $false, $true, 1, $false, "342×8=", 2) | Out-Null
$d.Content.Find.Execute("154×4=", $true, $false, $false, $false, $false, $true, 1, $false, "733×4=", 2) | Out-Null
$d.Content.Find.Execute("877×6=", $true, $false, $false, $false, $false, $true, 1, $false, "112×9=", 2) | Out-Null
$d.Content.Find.Execute("649×6=", $true, $false, $false, $false, $false, $true, 1, $false, "174×9=", 2) | Out-Null
$d.Content.Find.Execute("796×5=", $true, $false, $false, $false, $false, $true, 1, $false, "237×6=", 2) | Out-Null
$d.Content.Find.Execute("355×5=", $true, $false, $false, $false, $false, $true, 1, $false, "380×3=", 2) | Out-Null
$d.Content.Find.Execute("203×9=", $true, $false, $false, $false, $false, $true, 1, $false, "714×9=", 2) | Out-Null
$d.Content.Find.Execute("602×9=", $true, $false, $false, $false, $false, $true, 1, $false, "133×9=", 2) | Out-Null
$d.Content.Find.Execute("122×2=", $true, $false, $false, $false, $false, $true, 1, $false, "687×5=", 2) | Out-Null
$d.Content.Find.Execute("241×5=", $true, $false, $false, $false, $false, $true, 1, $false, "554×2=", 2) | Out-Null
$d.Content.Find.Execute("738×9=", $true, $false, $false, $false, $false, $true, 1, $false, "724×5=", 2) | Out-Null
